# Checkin Week 14 newsletter.
# Update the Rockies Injuries worksheet:
#  - Ian Desmond's row gets a new Last.Updated date and new Injury.Details text
#  - Carlos Gonzalez's and Gerardo Parra's rows are removed (no longer injured / updated)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Ian Desmond row (row 5): Last.Updated (col C) and Injury.Details (col E)
$ws.Cells.Item(5, 3).Value = "July 03 2017"
$ws.Cells.Item(5, 5).Value = "Desmond has landed on the 10-day disabled list with a right calf strain and it is unknown as to how long he will be out of action."

# Remove the Carlos Gonzalez (row 6) and Gerardo Parra (row 7) rows entirely
$ws.Rows.Item(6).Resize(2).Delete()

# Move the active selection as recorded in the saved workbook
$ws.Range("A28").Select()
